# Auto-generated edit script: updates cryptos list values (prices & volume %)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.604.41'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '1.811.45'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("E6").Value = '  +3.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '36.61'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.294'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0685'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0968'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").Value = '2.072.94'
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.34'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.06%  '
$ws.Range("D14").Value = '1.819.67'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.633'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '34.547.99'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.89%  '
$ws.Range("D20").Value = '0.0₃0779'
$ws.Range("E20").Value = '  -2.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("E24").Value = '  +4.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("E28").Value = '  +1.89%  '
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.93'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0517'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.43%  '
$ws.Range("D35").Value = '1.364.10'
$ws.Range("E35").Value = '  -2.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.655'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.73%  '
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("E38").Value = '  -5.11%  '
$ws.Range("E39").Value = '  -1.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.43'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '81.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.51%  '
$ws.Range("E42").Value = '  -1.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.939'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.65'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("D47").Value = '1.973.56'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("E49").Value = '  -0.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.79'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("D51").Value = '0.0₆0122'
$ws.Range("E51").Value = '  -6.46%  '
